# Update "想去人数" (want-to-go count) column F values on the "展览" and
# "全部类型" worksheets to reflect newly generated output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# Row => new F-column value, identical update set for both sheets
$commonUpdates = @{
    2  = 1064
    3  = 758
    4  = 261
    5  = 32
    8  = 1790
    9  = 6414
    10 = 486
    11 = 373
    12 = 314
    15 = 139
    16 = 6402
    17 = 276
    18 = 1292
    23 = 279
    27 = 99
    28 = 12
    29 = 393
    30 = 93
    33 = 49
}

foreach ($row in $commonUpdates.Keys) {
    $value = $commonUpdates[$row]
    $ws1.Cells.Item($row, 6).Value = $value
    $ws4.Cells.Item($row, 6).Value = $value
}

# Sheet-specific updates (different before-values on each sheet)
$ws1.Cells.Item(13, 6).Value = 103
$ws1.Cells.Item(24, 6).Value = 110

$ws4.Cells.Item(13, 6).Value = 103
$ws4.Cells.Item(24, 6).Value = 110
$ws4.Cells.Item(37, 6).Value = 63
